# Applies the "Sync attendance_reports, modules_schedules, and assets from
# main repo - 2026-01-06 09:21:29" update to the session-analysis sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Every "Recorded By" cell that lists "System, dnasr281@gmail.com"
#    now lists the same two recorders in the opposite order.
# ---------------------------------------------------------------------
$used = $ws.UsedRange
$used.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System", -4142, 1, $false) | Out-Null

# ---------------------------------------------------------------------
# 2. Overall statistics block (K4:L10) - Missing / Pending session totals.
# ---------------------------------------------------------------------
$ws.Range("L7").Value = 21   # Missing Sessions
$ws.Range("L8").Value = 84   # Pending Sessions

# ---------------------------------------------------------------------
# 3. Per-group statistics (rows 21-26 = B1D1, B1D2, B1E1, B1E2, B1F1, B1F2)
#    one previously "Pending" session in each of these groups is now
#    counted as "Missing" instead.
# ---------------------------------------------------------------------
$ws.Range("P21").Value = 2
$ws.Range("Q21").Value = 8

$ws.Range("P22").Value = 2
$ws.Range("Q22").Value = 8

$ws.Range("P23").Value = 2
$ws.Range("Q23").Value = 8

$ws.Range("P24").Value = 3
$ws.Range("Q24").Value = 8

$ws.Range("P25").Value = 2
$ws.Range("Q25").Value = 8

$ws.Range("P26").Value = 2
$ws.Range("Q26").Value = 8

# ---------------------------------------------------------------------
# 4. The 06/01/2026 session for groups B1D1, B1D2, B1E1, B1E2, B1F1, B1F2
#    (rows 176, 203, 230, 257, 284, 311) flips from "Pending" to
#    "Not Recorded", which also changes its row shading from the
#    "Pending" yellow style to the "Not Recorded" pink style (copy the
#    formatting from row 3, an existing "Not Recorded" row).
# ---------------------------------------------------------------------
$statusRows = @(176, 203, 230, 257, 284, 311)

$formatSource = $ws.Range("A3:I3")
$formatSource.Copy() | Out-Null

foreach ($rowNum in $statusRows) {
    $target = $ws.Range("A" + $rowNum + ":I" + $rowNum)
    $target.PasteSpecial(-4122) | Out-Null
    $ws.Range("I" + $rowNum).Value = "Not Recorded"
}

$excel.CutCopyMode = 0
